# Apply the "error solve ifrs list" correction to the company_list sheet.
# Rows 2-6 (FY2014-FY2018 IFRS-consolidated figures) get corrected values;
# rows 7-9 (FY2019E-FY2021E estimate rows) had their figures removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param($ws, [int]$row, [hashtable]$values)
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}

# Row 2 (FY2014/12 IFRS-consolidated)
Set-RowValues $ws 2 @{
    "D" = 1851;  "E" = 107;    "F" = 107;    "G" = 111;    "H" = 90;
    "I" = 90;    "J" = 0;      "K" = 1560;   "L" = 799;    "M" = 761;
    "N" = 756;   "O" = 5;      "P" = 249;    "Q" = 169;    "R" = -125;
    "S" = -42;   "T" = 27;     "U" = 142;    "V" = 521;    "W" = 5.77;
    "X" = 4.89;  "Y" = 12.74;  "Z" = 5.94;   "AA" = 104.99; "AB" = 194.54;
    "AC" = 363;  "AD" = 5.35;  "AE" = 3134;  "AF" = 0.62;  "AG" = 45;
    "AH" = 2.32; "AI" = 12.01; "AJ" = 24939425;
}

# Row 3 (FY2015/12 IFRS-consolidated)
Set-RowValues $ws 3 @{
    "D" = 1854;  "E" = 68;     "F" = 68;     "G" = 72;     "H" = 54;
    "I" = 56;    "J" = -2;     "K" = 1600;   "L" = 806;    "M" = 793;
    "N" = 770;   "O" = 23;     "P" = 249;    "Q" = 96;     "R" = -74;
    "S" = -16;   "T" = 38;     "U" = 58;     "V" = 585;    "W" = 3.65;
    "X" = 2.91;  "Y" = 7.36;   "Z" = 3.42;   "AA" = 101.61; "AB" = 209.61;
    "AC" = 225;  "AD" = 9.119999999999999; "AE" = 3191; "AF" = 0.64; "AG" = 30;
    "AH" = 1.46; "AI" = 12.89; "AJ" = 24939425;
}

# Row 4 (FY2016/12 IFRS-consolidated)
Set-RowValues $ws 4 @{
    "D" = 1723;  "E" = 67;     "F" = 67;     "G" = 93;     "H" = 69;
    "I" = 70;    "J" = -1;     "K" = 1678;   "L" = 828;    "M" = 850;
    "N" = 834;   "O" = 16;     "P" = 249;    "Q" = 158;    "R" = -119;
    "S" = -35;   "T" = 37;     "U" = 121;    "V" = 564;    "W" = 3.88;
    "X" = 3.98;  "Y" = 8.66;   "Z" = 4.18;   "AA" = 97.42; "AB" = 233;
    "AC" = 279;  "AD" = 23.93; "AE" = 3456;  "AF" = 1.93;  "AG" = 30;
    "AH" = 0.45; "AI" = 10.42; "AJ" = 24939425;
}

# Row 5 (FY2017/12 IFRS-consolidated)
Set-RowValues $ws 5 @{
    "D" = 1745;  "E" = 59;     "F" = 59;     "G" = 79;     "H" = 51;
    "I" = 51;    "J" = 0;      "K" = 1831;   "L" = 944;    "M" = 888;
    "N" = 888;   "O" = 0;      "P" = 249;    "Q" = 10;     "R" = -107;
    "S" = 115;   "T" = 42;     "U" = -32;    "V" = 666;    "W" = 3.41;
    "X" = 2.9;   "Y" = 5.93;   "Z" = 2.89;   "AA" = 106.32; "AB" = 259.04;
    "AC" = 205;  "AD" = 13.13; "AE" = 3598;  "AF" = 0.75;  "AG" = 35;
    "AH" = 1.3;  "AI" = 16.9;  "AJ" = 24939425;
}

# Row 6 (FY2018/12 IFRS-consolidated) -- note: no J/O columns, same as original
Set-RowValues $ws 6 @{
    "D" = 1584;  "E" = 58;     "F" = 58;     "G" = 48;     "H" = 37;
    "I" = 39;    "K" = 2136;   "L" = 1212;   "M" = 924;    "N" = 926;
    "P" = 249;   "Q" = -57;    "R" = -185;   "S" = 262;    "T" = 151;
    "U" = -209;  "V" = 935;    "W" = 3.67;   "X" = 2.34;   "Y" = 4.34;
    "Z" = 1.87;  "AA" = 131.22; "AB" = 274.72; "AC" = 158; "AD" = 13.74;
    "AE" = 3754; "AF" = 0.58;  "AG" = 30;    "AH" = 1.38;  "AI" = 18.79;
    "AJ" = 24939425;
}

# Rows 7-9 (FY2019E/2020E/2021E estimate rows): figures removed entirely,
# leaving only the A (index), B (label) and C (period) columns.
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
